$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 211 (shifts existing rows 211-234 down to 212-235,
# carrying their formatting along - mirrors Excel's native row insert).
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new weekly price record.
$ws.Range("A211").Value = 11
$ws.Range("B211").Value = "Vega Monumental Concepción"
$ws.Range("C211").Value = "Bíobío"
$ws.Range("D211").Value = 45154
$ws.Range("E211").Value = 8
$ws.Range("F211").Value = 100112043
$ws.Range("G211").Value = "Pepino ensalada"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 220
$ws.Range("K211").Value = 9000
$ws.Range("L211").Value = 10000
$ws.Range("M211").Value = 9455
$ws.Range("N211").Value = "`$/caja 60 unidades"
$ws.Range("O211").Value = "Región de Arica y Parinacota"
$ws.Range("P211").Value = 158
$ws.Range("Q211").Value = 60
$ws.Range("R211").Value = "Hortaliza"
